# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interest-count) column F figures across all four
# sheets (展览 / 演出 / 本地生活 / 全部类型) to the freshly scraped counts.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 13800
$ws.Cells.Item(5, 6).Value = 76
$ws.Cells.Item(6, 6).Value = 802
$ws.Cells.Item(7, 6).Value = 2192
$ws.Cells.Item(8, 6).Value = 193
$ws.Cells.Item(9, 6).Value = 134
$ws.Cells.Item(10, 6).Value = 114
$ws.Cells.Item(11, 6).Value = 238
$ws.Cells.Item(13, 6).Value = 609
$ws.Cells.Item(14, 6).Value = 464
$ws.Cells.Item(15, 6).Value = 526
$ws.Cells.Item(16, 6).Value = 346
$ws.Cells.Item(17, 6).Value = 37
$ws.Cells.Item(18, 6).Value = 316
$ws.Cells.Item(19, 6).Value = 895
$ws.Cells.Item(20, 6).Value = 161
$ws.Cells.Item(21, 6).Value = 87
$ws.Cells.Item(22, 6).Value = 39
$ws.Cells.Item(25, 6).Value = 108
$ws.Cells.Item(26, 6).Value = 45
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 56
$ws.Cells.Item(6, 6).Value = 130
$ws.Cells.Item(8, 6).Value = 2052
$ws.Cells.Item(13, 6).Value = 88
$ws.Cells.Item(15, 6).Value = 1886
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 237
$ws.Cells.Item(3, 6).Value = 189
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 237
$ws.Cells.Item(3, 6).Value = 13800
$ws.Cells.Item(6, 6).Value = 76
$ws.Cells.Item(7, 6).Value = 802
$ws.Cells.Item(8, 6).Value = 56
$ws.Cells.Item(10, 6).Value = 2192
$ws.Cells.Item(11, 6).Value = 189
$ws.Cells.Item(12, 6).Value = 193
$ws.Cells.Item(13, 6).Value = 134
$ws.Cells.Item(14, 6).Value = 114
$ws.Cells.Item(15, 6).Value = 238
$ws.Cells.Item(19, 6).Value = 130
$ws.Cells.Item(21, 6).Value = 609
$ws.Cells.Item(22, 6).Value = 464
$ws.Cells.Item(23, 6).Value = 526
$ws.Cells.Item(24, 6).Value = 346
$ws.Cells.Item(25, 6).Value = 37
$ws.Cells.Item(26, 6).Value = 316
$ws.Cells.Item(27, 6).Value = 895
$ws.Cells.Item(29, 6).Value = 2052
$ws.Cells.Item(34, 6).Value = 161
$ws.Cells.Item(35, 6).Value = 87
$ws.Cells.Item(36, 6).Value = 39
$ws.Cells.Item(38, 6).Value = 88
$ws.Cells.Item(41, 6).Value = 108
$ws.Cells.Item(42, 6).Value = 45
$ws.Cells.Item(43, 6).Value = 1886
